$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New transaction rows for week ending 2021-04-04 (54 rows: 506-559)
# Columns: A=Datum, B=Receipt Number, C=Konto, D=Beskrivning, E=Debet, F=Kredit
$rows = @(
    @(44284, "Reko61", 3011, "Reko Swish +46703533270", "", 460.71),
    @(44284, "Reko61", 2611, "Reko Swish +46703533270", "", 55.29),
    @(44284, "Reko61", 1930, "Reko Swish +46703533270", 516, ""),
    @(44284, "Reko62", 3011, "Reko Swish +46767736486", "", 371.43),
    @(44284, "Reko62", 2611, "Reko Swish +46767736486", "", 44.57),
    @(44284, "Reko62", 1930, "Reko Swish +46767736486", 416, ""),
    @(44284, "Reko63", 3011, "Reko Swish +46721459019", "", 691.0700000000001),
    @(44284, "Reko63", 2611, "Reko Swish +46721459019", "", 82.93000000000001),
    @(44284, "Reko63", 1930, "Reko Swish +46721459019", 774, ""),
    @(44284, "Reko64", 3011, "Reko Swish +46737032257", "", 230.36),
    @(44284, "Reko64", 2611, "Reko Swish +46737032257", "", 27.64),
    @(44284, "Reko64", 1930, "Reko Swish +46737032257", 258, ""),
    @(44284, "Reko65", 3011, "Reko Swish +46702331968", "", 616.0700000000001),
    @(44284, "Reko65", 2611, "Reko Swish +46702331968", "", 73.93000000000001),
    @(44284, "Reko65", 1930, "Reko Swish +46702331968", 690, ""),
    @(44284, "Reko66", 3011, "Reko Swish +46709941173", "", 230.36),
    @(44284, "Reko66", 2611, "Reko Swish +46709941173", "", 27.64),
    @(44284, "Reko66", 1930, "Reko Swish +46709941173", 258, ""),
    @(44284, "Reko67", 3011, "Reko Swish +46737600861", "", 230.36),
    @(44284, "Reko67", 2611, "Reko Swish +46737600861", "", 27.64),
    @(44284, "Reko67", 1930, "Reko Swish +46737600861", 258, ""),
    @(44284, "", 4010, "M&S RB BROMMA K0135", 84.90000000000001, ""),
    @(44284, "", 2645, "M&S RB BROMMA K0135", 10.19, ""),
    @(44284, "", 1930, "M&S RB BROMMA K0135", "", 95.09),
    @(44284, "", 4010, "M&S RB BROMMA K0135", 328.1, ""),
    @(44284, "", 2645, "M&S RB BROMMA K0135", 39.37, ""),
    @(44284, "", 1930, "M&S RB BROMMA K0135", "", 367.47),
    @(44285, "Reko68", 3011, "Reko Swish +46703384055", "", 230.36),
    @(44285, "Reko68", 2611, "Reko Swish +46703384055", "", 27.64),
    @(44285, "Reko68", 1930, "Reko Swish +46703384055", 258, ""),
    @(44285, "Reko69", 3011, "Reko Swish +46708649109", "", 616.0700000000001),
    @(44285, "Reko69", 2611, "Reko Swish +46708649109", "", 73.93000000000001),
    @(44285, "Reko69", 1930, "Reko Swish +46708649109", 690, ""),
    @(44285, "Reko70", 3011, "Reko Swish +46704545110", "", 460.71),
    @(44285, "Reko70", 2611, "Reko Swish +46704545110", "", 55.29),
    @(44285, "Reko70", 1930, "Reko Swish +46704545110", 516, ""),
    @(44285, "Reko71", 3011, "Reko Swish +46723047499", "", 423.21),
    @(44285, "Reko71", 2611, "Reko Swish +46723047499", "", 50.79),
    @(44285, "Reko71", 1930, "Reko Swish +46723047499", 474, ""),
    @(44285, "Reko72", 3011, "Reko Swish +46703242407", "", 616.0700000000001),
    @(44285, "Reko72", 2611, "Reko Swish +46703242407", "", 73.93000000000001),
    @(44285, "Reko72", 1930, "Reko Swish +46703242407", 690, ""),
    @(44285, "Reko73", 3011, "Reko Swish +46707964655", "", 211.61),
    @(44285, "Reko73", 2611, "Reko Swish +46707964655", "", 25.39),
    @(44285, "Reko73", 1930, "Reko Swish +46707964655", 237, ""),
    @(44286, "Reko74", 3011, "Reko Swish +46704608939", "", 230.36),
    @(44286, "Reko74", 2611, "Reko Swish +46704608939", "", 27.64),
    @(44286, "Reko74", 1930, "Reko Swish +46704608939", 258, ""),
    @(44287, "1012353", 3011, "Order 1012353 Card(Stripe)", "", 1062.5),
    @(44287, "1012353", 2611, "Order 1012353 Card(Stripe)", "", 127.5),
    @(44287, "1012353", 1930, "Order 1012353 Card(Stripe)", 1190, ""),
    @(44290, "", 4010, "NGROCERIES K0135", 172.32, ""),
    @(44290, "", 2645, "NGROCERIES K0135", 20.68, ""),
    @(44290, "", 1930, "NGROCERIES K0135", "", 193)
)

$startRow = 506
$dateFormat = "YYYY-MM-DD HH:MM:SS"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $datum = $data[0]
    $receiptNumber = $data[1]
    $konto = $data[2]
    $beskrivning = $data[3]
    $debet = $data[4]
    $kredit = $data[5]

    # A: Datum - date-formatted like the rest of the column (style matches existing A2:A505)
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = $dateFormat
    $cellA.Value = $datum

    # B: Receipt Number (always stored as text, even when it looks numeric,
    # so values like "1012353" aren't silently re-typed as a number)
    if ($receiptNumber -ne "") {
        $cellB = $ws.Cells.Item($r, 2)
        $cellB.NumberFormat = "@"
        $cellB.Value = $receiptNumber
        $cellB.Style = "Normal"
    }

    # C: Konto (plain number)
    $ws.Cells.Item($r, 3).Value = $konto

    # D: Beskrivning (text)
    $cellD = $ws.Cells.Item($r, 4)
    $cellD.NumberFormat = "@"
    $cellD.Value = $beskrivning
    $cellD.Style = "Normal"

    # E: Debet (only one of Debet/Kredit is populated per row)
    if ($debet -ne "") {
        $ws.Cells.Item($r, 5).Value = $debet
    }

    # F: Kredit
    if ($kredit -ne "") {
        $ws.Cells.Item($r, 6).Value = $kredit
    }
}
